# import_account_opening: new features
# - add a second "Prima Alpha S.p.A." customer line (new row 3)
# - remove the VAT number from the "Freie Universität Berlin" line
# - add a new "Mario Rossi" supplier line with an Italian fiscal code (new row 7)
# - widen column E a bit and move the active selection to E6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new rows first, from top to bottom, so row numbers below
#     refer to the already-shifted sheet. ---

# New row 3: duplicate "Prima Alpha S.p.A." entry (old rows 3-6 shift to 4-7)
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 right away (Nome / Fornitore / Partita IVA / Avere)
$ws.Range("B3").Value = "Prima Alpha S.p.A."
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "IT00115719999"
$ws.Range("H3").Value = 250

# New row 7: "Mario Rossi" entry (old row7 "Banca" shifts to row 8)
$ws.Rows.Item(7).Insert()

# Fill in the new row 7 (Nome / Cliente / Partita IVA / Dare)
$ws.Range("B7").Value = "Mario Rossi"
$ws.Range("C7").Value = 1
$ws.Range("E7").Value = "RSSMRA69C02D612M"
$ws.Range("G7").Value = 50

# Remove the VAT number previously on the "Freie Universität Berlin" row (now row 6)
$ws.Range("E6").Value = ""

# Widen column E slightly
$ws.Columns.Item(5).ColumnWidth = 18.6

# Move the active selection to E6, matching the saved view state
$ws.Range("E6").Select() | Out-Null
